# Applies the cryptos list update described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.575.76"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").Value = "2.290.16"
$ws.Range("E3").Value = "  -1.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.23"
$ws.Range("E5").Value = "  -1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.00"
$ws.Range("E6").Value = "  -5.38%  "

$ws.Range("E7").Value = "  -1.93%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -3.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.37"
$ws.Range("E10").Value = "  -5.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.94"
$ws.Range("E12").Value = "  -5.10%  "

$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.59"
$ws.Range("E14").Value = "  +5.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.74"
$ws.Range("E15").Value = "  -0.72%  "

$ws.Range("D16").Value = "2.645.27"
$ws.Range("E16").Value = "  -1.86%  "

$ws.Range("D17").Value = "2.288.89"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.795"
$ws.Range("E18").Value = "  -1.18%  "

$ws.Range("D19").Value = "42.467.12"
$ws.Range("E19").Value = "  -2.11%  "

$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.47"
$ws.Range("E21").Value = "  -3.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.99"
$ws.Range("E22").Value = "  -1.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.82"
$ws.Range("E23").Value = "  -2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.91"
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("E27").Value = "  -3.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.34"
$ws.Range("E28").Value = "  -2.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.75"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.75"
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.11"
$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("E34").Value = "  +3.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.95"
$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.00"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0691"
$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.81"
$ws.Range("E39").Value = "  -3.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0996"
$ws.Range("E40").Value = "  -2.60%  "

$ws.Range("E41").Value = "  -4.66%  "

$ws.Range("E42").Value = "  -1.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.36"
$ws.Range("E43").Value = "  -6.84%  "

$ws.Range("D44").Value = "1.958.89"
$ws.Range("E44").Value = "  -0.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0279"
$ws.Range("E45").Value = "  -1.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.75"
$ws.Range("E46").Value = "  -3.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.47"
$ws.Range("E47").Value = "  -6.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.81"
$ws.Range("E48").Value = "  -3.65%  "

$ws.Range("D49").Value = "2.513.75"
$ws.Range("E49").Value = "  -1.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.60"
$ws.Range("E50").Value = "  -6.03%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.74"
$ws.Range("E51").Value = "  -2.16%  "
